# Commit: "Update Excel SCD0011 until SCD0016"
#
# The source sheet (originally cloned from a "SCD0187" template) is being
# renumbered to SCD0011, and its TC_ID cell is updated to match the new
# scenario id. The now-unused old TC_ID text ("DGS-202") is dropped from
# the shared-string table automatically when no cell references it anymore.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet tab: SCD0187 -> SCD0011
$ws.Name = "SCD0011"

# TC_ID cell (B2) now holds the new scenario id
$ws.Range("B2").Value = "SCD0011-018"

# Cell B3 is the active selection on this sheet
$ws.Range("B3").Select() | Out-Null

# Column B needs to be a bit wider to fit the longer "SCD0011-018" id
$ws.Columns.Item(2).ColumnWidth = 11.67
